# Update the "Förändrad" (Changed) date column (C2:C6) from 2023-09-15 to
# 2023-09-16 (serial date 45184 -> 45185), keeping existing number formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$newDate = [DateTime]::FromOADate(45185)

for ($row = 2; $row -le 6; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
